$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new columns: P (Platform Account) and Q (Connect Account) ---
$ws.Range("P1").Value = "Platform Account"
$ws.Range("Q1").Value = "Connect Account"

$ws.Range("P2").Value = "EatMe - POS 2 (acct_1O9kwUAaoVAZ6m8M)"
$ws.Range("Q2").Value = "acct_1REm1MPMERGGWtpY"

# Column widths for the new columns (COM ColumnWidth is offset by 5/6 from the
# raw OOXML <col width>, so subtract that to land on the exact target widths)
$ws.Columns.Item(16).ColumnWidth = 32 - (5/6)
$ws.Columns.Item(17).ColumnWidth = 25 - (5/6)

# --- Update the selection shown in the saved view ---
[void]$ws.Range("O20").Select()
